# PipeDiameters.xlsx -- "updated file naming convention in logger"
#
# Extend the second lookup table (Gas Demand vs Pipe Diameter) with two more
# rows (2800, 2900 kg/s plant loads), keep the shared B-column formula
# going, nudge a few column widths, and move the on-screen view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows 33 & 34 ---------------------------------------------------
$ws.Range("A33").Value = 2800
$ws.Range("A34").Value = 2900

# Match the look of the existing A22:A32 cells (fill + border) by copying
# their format onto A32 (which was missing it) and the two new cells.
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A32:A34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Formula for the new B cells follows the same pattern as B21:B32
# ( = (A<row> * 1,000,000) / ($E$20 * 1000 * $E$21) ), giving
# 411.54985265045457 and 426.24806167368507 respectively.
$ws.Range("B33").Formula = '=(A33*1000000)/($E$20*1000*$E$21)'
$ws.Range("B34").Formula = '=(A34*1000000)/($E$20*1000*$E$21)'

# Carry the B-column formatting (shaded fill + border + number format) down
# onto the two new formula cells too.
$ws.Range("B32").Copy() | Out-Null
$ws.Range("B33:B34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column width tweaks -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 28.33   # -> ~29.16 chars (was 43.5)
$ws.Columns.Item(3).ColumnWidth = 14.65   # -> 15.5 chars   (was 22.66)
$ws.Columns.Item(4).ColumnWidth = 22.33   # -> ~23.16 chars (was 29.66)

# --- View / selection state ---------------------------------------------------
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

$wb.Windows.Item(1).WindowState = -4143   # xlNormal (un-minimize)
